$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E width adjustment (merged range in target; width quantizes to nearest
# achievable pixel-based value through the Excel object model)
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 9.75

# Populate new column R (2021 data) - copy formatting from column Q, then set value
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 2021
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 297.10000000000002
$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 311
$ws.Range("Q7").Copy($ws.Range("R7"))
$ws.Range("R7").Value = 283.3
$ws.Range("Q8").Copy($ws.Range("R8"))
$ws.Range("R8").Value = 281.7
$ws.Range("Q9").Copy($ws.Range("R9"))
$ws.Range("R9").Value = 299.5
$ws.Range("Q10").Copy($ws.Range("R10"))
$ws.Range("R10").Value = 263.3
$ws.Range("Q11").Copy($ws.Range("R11"))
$ws.Range("R11").Value = 287.10000000000002
$ws.Range("Q12").Copy($ws.Range("R12"))
$ws.Range("R12").Value = 313.7
$ws.Range("Q13").Copy($ws.Range("R13"))
$ws.Range("R13").Value = 260.10000000000002
$ws.Range("Q14").Copy($ws.Range("R14"))
$ws.Range("R14").Value = 358.2
$ws.Range("Q15").Copy($ws.Range("R15"))
$ws.Range("R15").Value = 360.7
$ws.Range("Q16").Copy($ws.Range("R16"))
$ws.Range("R16").Value = 355.7
$ws.Range("Q17").Copy($ws.Range("R17"))
$ws.Range("R17").Value = 332.2
$ws.Range("Q18").Copy($ws.Range("R18"))
$ws.Range("R18").Value = 377.1
$ws.Range("Q19").Copy($ws.Range("R19"))
$ws.Range("R19").Value = 285.60000000000002
$ws.Range("Q20").Copy($ws.Range("R20"))
$ws.Range("R20").Value = 248.1
$ws.Range("Q21").Copy($ws.Range("R21"))
$ws.Range("R21").Value = 259.2
$ws.Range("Q22").Copy($ws.Range("R22"))
$ws.Range("R22").Value = 236.8
$ws.Range("Q23").Copy($ws.Range("R23"))
$ws.Range("R23").Value = 347.9
$ws.Range("Q24").Copy($ws.Range("R24"))
$ws.Range("R24").Value = 361.7
$ws.Range("Q25").Copy($ws.Range("R25"))
$ws.Range("R25").Value = 333.8
$ws.Range("Q26").Copy($ws.Range("R26"))
$ws.Range("R26").Value = 391.1
$ws.Range("Q27").Copy($ws.Range("R27"))
$ws.Range("R27").Value = 400.7
$ws.Range("Q28").Copy($ws.Range("R28"))
$ws.Range("R28").Value = 381.9
$ws.Range("Q29").Copy($ws.Range("R29"))
$ws.Range("R29").Value = 255.7
$ws.Range("Q30").Copy($ws.Range("R30"))
$ws.Range("R30").Value = 256
$ws.Range("Q31").Copy($ws.Range("R31"))
$ws.Range("R31").Value = 255.5
$ws.Range("Q32").Copy($ws.Range("R32"))
$ws.Range("R32").Value = 258.89999999999998
$ws.Range("Q33").Copy($ws.Range("R33"))
$ws.Range("R33").Value = 271.5
$ws.Range("Q34").Copy($ws.Range("R34"))
$ws.Range("R34").Value = 247
$ws.Range("Q35").Copy($ws.Range("R35"))
$ws.Range("Q36").Copy($ws.Range("R36"))
$ws.Range("R36").Value = 59.5
$ws.Range("Q37").Copy($ws.Range("R37"))
$ws.Range("R37").Value = 62.7
$ws.Range("Q38").Copy($ws.Range("R38"))
$ws.Range("R38").Value = 56.4
$ws.Range("Q39").Copy($ws.Range("R39"))
$ws.Range("R39").Value = 37.799999999999997
$ws.Range("Q40").Copy($ws.Range("R40"))
$ws.Range("R40").Value = 31.9
$ws.Range("Q41").Copy($ws.Range("R41"))
$ws.Range("R41").Value = 43.8
$ws.Range("Q42").Copy($ws.Range("R42"))
$ws.Range("R42").Value = 46.2
$ws.Range("Q43").Copy($ws.Range("R43"))
$ws.Range("R43").Value = 48.9
$ws.Range("Q44").Copy($ws.Range("R44"))
$ws.Range("R44").Value = 43.6
$ws.Range("Q45").Copy($ws.Range("R45"))
$ws.Range("R45").Value = 82.4
$ws.Range("Q46").Copy($ws.Range("R46"))
$ws.Range("R46").Value = 93.3
$ws.Range("Q47").Copy($ws.Range("R47"))
$ws.Range("R47").Value = 71.5
$ws.Range("Q48").Copy($ws.Range("R48"))
$ws.Range("R48").Value = 71.599999999999994
$ws.Range("Q49").Copy($ws.Range("R49"))
$ws.Range("R49").Value = 81.7
$ws.Range("Q50").Copy($ws.Range("R50"))
$ws.Range("R50").Value = 61.2
$ws.Range("Q51").Copy($ws.Range("R51"))
$ws.Range("R51").Value = 45.9
$ws.Range("Q52").Copy($ws.Range("R52"))
$ws.Range("R52").Value = 50.3
$ws.Range("Q53").Copy($ws.Range("R53"))
$ws.Range("R53").Value = 41.4
$ws.Range("Q54").Copy($ws.Range("R54"))
$ws.Range("R54").Value = 52.1
$ws.Range("Q55").Copy($ws.Range("R55"))
$ws.Range("R55").Value = 58.8
$ws.Range("Q56").Copy($ws.Range("R56"))
$ws.Range("R56").Value = 45.3
$ws.Range("Q57").Copy($ws.Range("R57"))
$ws.Range("R57").Value = 94.9
$ws.Range("Q58").Copy($ws.Range("R58"))
$ws.Range("R58").Value = 99.5
$ws.Range("Q59").Copy($ws.Range("R59"))
$ws.Range("R59").Value = 90.4
$ws.Range("Q60").Copy($ws.Range("R60"))
$ws.Range("R60").Value = 61.4
$ws.Range("Q61").Copy($ws.Range("R61"))
$ws.Range("R61").Value = 62.7
$ws.Range("Q62").Copy($ws.Range("R62"))
$ws.Range("R62").Value = 60.2
$ws.Range("Q63").Copy($ws.Range("R63"))
$ws.Range("R63").Value = 54.9
$ws.Range("Q64").Copy($ws.Range("R64"))
$ws.Range("R64").Value = 54.5
$ws.Range("Q65").Copy($ws.Range("R65"))
$ws.Range("R65").Value = 55.2
$ws.Range("Q66").Copy($ws.Range("R66"))
$ws.Range("Q67").Copy($ws.Range("R67"))
$ws.Range("R67").Value = 7.1
$ws.Range("Q68").Copy($ws.Range("R68"))
$ws.Range("R68").Value = 5.8
$ws.Range("Q69").Copy($ws.Range("R69"))
$ws.Range("R69").Value = 8.3000000000000007
$ws.Range("Q70").Copy($ws.Range("R70"))
$ws.Range("R70").Value = 10.8
$ws.Range("Q71").Copy($ws.Range("R71"))
$ws.Range("R71").Value = 6.4
$ws.Range("Q72").Copy($ws.Range("R72"))
$ws.Range("R72").Value = 15.5
$ws.Range("Q73").Copy($ws.Range("R73"))
$ws.Range("R73").Value = 13.1
$ws.Range("Q74").Copy($ws.Range("R74"))
$ws.Range("R74").Value = 10
$ws.Range("Q75").Copy($ws.Range("R75"))
$ws.Range("R75").Value = 16.3
$ws.Range("Q76").Copy($ws.Range("R76"))
$ws.Range("R76").Value = 6.9
$ws.Range("Q77").Copy($ws.Range("R77"))
$ws.Range("R77").Value = 5.2
$ws.Range("Q78").Copy($ws.Range("R78"))
$ws.Range("R78").Value = 8.6999999999999993
$ws.Range("Q79").Copy($ws.Range("R79"))
$ws.Range("R79").Value = 7.5
$ws.Range("Q80").Copy($ws.Range("R80"))
$ws.Range("R80").Value = 10.7
$ws.Range("Q81").Copy($ws.Range("R81"))
$ws.Range("R81").Value = 4.2
$ws.Range("Q82").Copy($ws.Range("R82"))
$ws.Range("R82").Value = 5.0999999999999996
$ws.Range("Q83").Copy($ws.Range("R83"))
$ws.Range("R83").Value = 3.7
$ws.Range("Q84").Copy($ws.Range("R84"))
$ws.Range("R84").Value = 6.6
$ws.Range("Q85").Copy($ws.Range("R85"))
$ws.Range("R85").Value = 3.7
$ws.Range("Q86").Copy($ws.Range("R86"))
$ws.Range("R86").Value = 4.4000000000000004
$ws.Range("Q87").Copy($ws.Range("R87"))
$ws.Range("R87").Value = 3
$ws.Range("Q88").Copy($ws.Range("R88"))
$ws.Range("R88").Value = 5.3
$ws.Range("Q89").Copy($ws.Range("R89"))
$ws.Range("R89").Value = 5.8
$ws.Range("Q90").Copy($ws.Range("R90"))
$ws.Range("R90").Value = 4.8
$ws.Range("Q91").Copy($ws.Range("R91"))
$ws.Range("R91").Value = 2.2000000000000002
$ws.Range("Q92").Copy($ws.Range("R92"))
$ws.Range("R92").Value = 2.5
$ws.Range("Q93").Copy($ws.Range("R93"))
$ws.Range("R93").Value = 1.9
$ws.Range("Q94").Copy($ws.Range("R94"))
$ws.Range("R94").Value = 9.5
$ws.Range("Q95").Copy($ws.Range("R95"))
$ws.Range("R95").Value = 5
$ws.Range("Q96").Copy($ws.Range("R96"))
$ws.Range("R96").Value = 13.7
$ws.Range("Q97").Copy($ws.Range("R97"))
$ws.Range("Q98").Copy($ws.Range("R98"))
$ws.Range("R98").Value = 14.2
$ws.Range("Q99").Copy($ws.Range("R99"))
$ws.Range("R99").Value = 17.3
$ws.Range("Q100").Copy($ws.Range("R100"))
$ws.Range("R100").Value = 11.2
$ws.Range("Q101").Copy($ws.Range("R101"))
$ws.Range("R101").Value = 7.4
$ws.Range("Q102").Copy($ws.Range("R102"))
$ws.Range("R102").Value = 6
$ws.Range("Q103").Copy($ws.Range("R103"))
$ws.Range("R103").Value = 8.8000000000000007
$ws.Range("Q104").Copy($ws.Range("R104"))
$ws.Range("R104").Value = 15
$ws.Range("Q105").Copy($ws.Range("R105"))
$ws.Range("R105").Value = 19.399999999999999
$ws.Range("Q106").Copy($ws.Range("R106"))
$ws.Range("R106").Value = 10.6
$ws.Range("Q107").Copy($ws.Range("R107"))
$ws.Range("R107").Value = 33.9
$ws.Range("Q108").Copy($ws.Range("R108"))
$ws.Range("R108").Value = 45.8
$ws.Range("Q109").Copy($ws.Range("R109"))
$ws.Range("R109").Value = 22.1
$ws.Range("Q110").Copy($ws.Range("R110"))
$ws.Range("R110").Value = 43.7
$ws.Range("Q111").Copy($ws.Range("R111"))
$ws.Range("R111").Value = 47.6
$ws.Range("Q112").Copy($ws.Range("R112"))
$ws.Range("R112").Value = 39.6
$ws.Range("Q113").Copy($ws.Range("R113"))
$ws.Range("R113").Value = 18.7
$ws.Range("Q114").Copy($ws.Range("R114"))
$ws.Range("R114").Value = 19.899999999999999
$ws.Range("Q115").Copy($ws.Range("R115"))
$ws.Range("R115").Value = 17.600000000000001
$ws.Range("Q116").Copy($ws.Range("R116"))
$ws.Range("R116").Value = 7.7
$ws.Range("Q117").Copy($ws.Range("R117"))
$ws.Range("R117").Value = 8.6999999999999993
$ws.Range("Q118").Copy($ws.Range("R118"))
$ws.Range("R118").Value = 6.7
$ws.Range("Q119").Copy($ws.Range("R119"))
$ws.Range("R119").Value = 6.7
$ws.Range("Q120").Copy($ws.Range("R120"))
$ws.Range("R120").Value = 10
$ws.Range("Q121").Copy($ws.Range("R121"))
$ws.Range("R121").Value = 3.6
$ws.Range("Q122").Copy($ws.Range("R122"))
$ws.Range("R122").Value = 3.7
$ws.Range("Q123").Copy($ws.Range("R123"))
$ws.Range("R123").Value = 5.5
$ws.Range("Q124").Copy($ws.Range("R124"))
$ws.Range("R124").Value = 2.1
$ws.Range("Q125").Copy($ws.Range("R125"))
$ws.Range("R125").Value = 9.8000000000000007
$ws.Range("Q126").Copy($ws.Range("R126"))
$ws.Range("R126").Value = 11.3
$ws.Range("Q127").Copy($ws.Range("R127"))
$ws.Range("R127").Value = 8.3000000000000007

# Set active selection to Q1 as in target sheetView
$ws.Range("Q1").Select()
